$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting the existing row 2 (and below) down to row 3
$ws.Rows.Item(2).Insert()

# New row 2 values
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 0.801303488973446
$ws.Range("C2").Value = 5.23404690932392
$ws.Range("D2").Value = 0.846169544288195

# Row 3 (previously row 2, shifted down) gets its own new values
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 0.657822406125166
$ws.Range("C3").Value = 6.99923359854542
$ws.Range("D3").Value = 0.731450191729779

# New row 4: A4=3, with B4/C4/D4 present but empty (no value)
$ws.Range("A4").Value = 3
$ws.Range("Z100").Copy()
$ws.Range("B4:D4").PasteSpecial(-4122)
